$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '43.868.34', '  -0.48%  ')
    ,@(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '2.351.24', '  -0.21%  ')
    ,@(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.00', '  -0.04%  ')
    ,@(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '239.99', '  +0.09%  ')
    ,@(6, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.668', '  -2.18%  ')
    ,@(7, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '74.07', '  -1.15%  ')
    ,@(8, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.00', '  -0.11%  ')
    ,@(9, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.598', '  -1.03%  ')
    ,@(10, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.101', '  +0.21%  ')
    ,@(11, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '61.86', '  +8.19%  ')
    ,@(12, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '33.28', '  +2.89%  ')
    ,@(13, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.108', '  +0.34%  ')
    ,@(14, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '7.25', '  -0.80%  ')
    ,@(15, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '16.21', '  -1.96%  ')
    ,@(16, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.903', '  -0.73%  ')
    ,@(17, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '2.348.18', '  -0.80%  ')
    ,@(18, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '43.803.91', '  -0.42%  ')
    ,@(19, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0000103', '  +0.19%  ')
    ,@(20, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '6.63', '  -1.16%  ')
    ,@(21, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '77.94', '  +1.02%  ')
    ,@(22, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '252.71', '  -1.77%  ')
    ,@(23, 'WEMIXToken', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '3.82', '  +2.54%  ')
    ,@(24, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.00', '  +0.06%  ')
    ,@(25, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '1.85', '  -0.79%  ')
    ,@(26, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '2.49', '  -0.79%  ')
    ,@(27, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '10.41', '  -3.21%  ')
    ,@(28, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.25', '  +0.11%  ')
    ,@(29, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '175.70', '  +0.46%  ')
    ,@(30, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '22.25', '  -2.82%  ')
    ,@(31, 'Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.127', '  -0.22%  ')
    ,@(32, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.135', '  -1.63%  ')
    ,@(33, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.0744', '  -2.16%  ')
    ,@(34, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '5.38', '  +0.07%  ')
    ,@(35, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '5.06', '  -4.64%  ')
    ,@(36, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '3.79', '  +1.55%  ')
    ,@(37, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.42', '  +2.00%  ')
    ,@(38, 'THORChain', 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune', '6.45', '  +1.41%  ')
    ,@(39, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.0272', '  -3.67%  ')
    ,@(40, 'FTXToken', 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt', '5.52', '  +15.60%  ')
    ,@(41, 'MultiversX', 'https://coinranking.com/coin/omwkOTglq+multiversx-egld', '65.49', '  +15.73%  ')
    ,@(42, 'InjectiveProtocol', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj', '19.29', '  -0.95%  ')
    ,@(43, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '9.15', '  +0.66%  ')
    ,@(44, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.201', '  -4.29%  ')
    ,@(45, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.105', '  -3.11%  ')
    ,@(46, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '0.999', '  -0.28%  ')
    ,@(47, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.23', '  -1.98%  ')
    ,@(48, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '2.44', '  -3.58%  ')
    ,@(49, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.15', '  -2.39%  ')
    ,@(50, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '98.29', '  -2.87%  ')
    ,@(51, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.90', '  +2.66%  ')
)

foreach ($r in $data) {
    $rowNum = $r[0]
    $arr = New-Object 'object[,]' 1,4
    $arr[0,0] = $r[1]
    $arr[0,1] = $r[2]
    $arr[0,2] = $r[3]
    $arr[0,3] = $r[4]
    $ws.Range("B$rowNum`:E$rowNum").Value = $arr
}
